$wb = $excel.ActiveWorkbook

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2281.1448
$ws.Range("I138").Value = 1782.8572
$ws.Range("J138").Value = 2471.4
$ws.Range("K138").Value = 5348.571599999999
$ws.Range("L138").Value = 7414.200000000001
$ws.Range("M138").Value = -208.5715999999993
$ws.Range("N138").Value = -17694.2

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26206.207
$ws.Range("I32").Value = 5084.587
$ws.Range("J32").Value = 107172.414
$ws.Range("K32").Value = 5084.587
$ws.Range("L32").Value = 107172.414
$ws.Range("M32").Value = -4797.587
$ws.Range("N32").Value = -107746.414

# ARM row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# ARM row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 24200
$ws.Range("J55").Value = 24200
$ws.Range("L55").Value = 24200
$ws.Range("N55").Value = -24830

# ARM row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# ARM row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19347
$ws.Range("I82").Value = 6419
$ws.Range("J82").Value = 38739
$ws.Range("K82").Value = 6419
$ws.Range("L82").Value = 38739
$ws.Range("M82").Value = -6036
$ws.Range("N82").Value = -39505

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 19347
$ws.Range("I85").Value = 6419
$ws.Range("J85").Value = 38739
$ws.Range("K85").Value = 6419
$ws.Range("L85").Value = 38739
$ws.Range("M85").Value = -5093
$ws.Range("N85").Value = -41391

# CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1429353.1
$ws.Range("I6").Value = 2501225.5
$ws.Range("J6").Value = 190
$ws.Range("K6").Value = 2501225.5
$ws.Range("L6").Value = 190
$ws.Range("M6").Value = -2501112.5
$ws.Range("N6").Value = -416

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32073.646
$ws.Range("I31").Value = 1886.069
$ws.Range("J31").Value = 207161.6
$ws.Range("K31").Value = 1886.069
$ws.Range("L31").Value = 207161.6
$ws.Range("M31").Value = -1591.069
$ws.Range("N31").Value = -207751.6

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 32073.646
$ws.Range("I34").Value = 1886.069
$ws.Range("J34").Value = 207161.6
$ws.Range("K34").Value = 1886.069
$ws.Range("L34").Value = 207161.6
$ws.Range("M34").Value = -1684.069
$ws.Range("N34").Value = -207565.6

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 29865
$ws.Range("J41").Value = 29865
$ws.Range("L41").Value = 29865
$ws.Range("N41").Value = -30721

# CRP row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 37999.332
$ws.Range("J51").Value = 37999.332
$ws.Range("L51").Value = 37999.332
$ws.Range("N51").Value = -39471.332

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 12421.8
$ws.Range("J60").Value = 14027.25
$ws.Range("L60").Value = 14027.25
$ws.Range("N60").Value = -15049.25

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 37999.332
$ws.Range("J61").Value = 37999.332
$ws.Range("L61").Value = 37999.332
$ws.Range("N61").Value = -38695.332

# CRP row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 36195
$ws.Range("J68").Value = 36195
$ws.Range("L68").Value = 36195
$ws.Range("N68").Value = -37693

# CRP row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 36195
$ws.Range("J71").Value = 36195
$ws.Range("L71").Value = 108585
$ws.Range("N71").Value = -116073

# CRP row 109
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 28995.8
$ws.Range("J109").Value = 28995.8
$ws.Range("L109").Value = 28995.8
$ws.Range("N109").Value = -31075.8

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 647.9
$ws.Range("I5").Value = 608.7778
$ws.Range("K5").Value = 1826.3334
$ws.Range("M5").Value = -1714.3334

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 745.2353000000001
$ws.Range("I86").Value = 568.875
$ws.Range("J86").Value = 902
$ws.Range("K86").Value = 1706.625
$ws.Range("L86").Value = 2706
$ws.Range("M86").Value = -520.625
$ws.Range("N86").Value = -5078

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 745.2353000000001
$ws.Range("I89").Value = 568.875
$ws.Range("J89").Value = 902
$ws.Range("K89").Value = 5119.875
$ws.Range("L89").Value = 8118
$ws.Range("M89").Value = 808.125
$ws.Range("N89").Value = -19974

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 647.9
$ws.Range("I135").Value = 608.7778
$ws.Range("K135").Value = 5479.000199999999
$ws.Range("M135").Value = -2944.000199999999

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 31274
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 31274
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 31274
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -32914

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3486.1177
$ws.Range("I122").Value = 3530.3333
$ws.Range("J122").Value = 3380
$ws.Range("K122").Value = 10590.9999
$ws.Range("L122").Value = 10140
$ws.Range("M122").Value = -8140.999899999999
$ws.Range("N122").Value = -15040

# GSM row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 15194.462
$ws.Range("J123").Value = 15194.462
$ws.Range("L123").Value = 15194.462
$ws.Range("N123").Value = -20094.462

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2014.2142
$ws.Range("I68").Value = 1300
$ws.Range("J68").Value = 2411
$ws.Range("K68").Value = 1300
$ws.Range("L68").Value = 2411
$ws.Range("M68").Value = -551
$ws.Range("N68").Value = -3909

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2014.2142
$ws.Range("I71").Value = 1300
$ws.Range("J71").Value = 2411
$ws.Range("K71").Value = 6500
$ws.Range("L71").Value = 12055
$ws.Range("M71").Value = -2756
$ws.Range("N71").Value = -19543

# LTW row 109
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 30516.666
$ws.Range("J109").Value = 30516.666
$ws.Range("L109").Value = 30516.666
$ws.Range("N109").Value = -33290.666

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2313
$ws.Range("I96").Value = 1750
$ws.Range("J96").Value = 2500.6667
$ws.Range("K96").Value = 1750
$ws.Range("L96").Value = 2500.6667
$ws.Range("M96").Value = -377
$ws.Range("N96").Value = -5246.6667
